$wb = $excel.ActiveWorkbook

# Report generated for handoff: set the "Latest Handoff Datetime" column (D)
# for the rows whose handoff file just got (re)generated, on both the
# zh-cn and de-de localization status sheets.

$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in 4,6,7,8,9,10) {
    $zhcn.Cells.Item($r, 4).Value = "2016-03-10 06:38:21"
}

$dede = $wb.Worksheets.Item("de-de")
foreach ($r in 4,6,7,8,9,10) {
    $dede.Cells.Item($r, 4).Value = "2016-03-10 06:38:30"
}
